# Update NATMI TPM-based ligand-receptor edge statistics (Sema4d-Erbb2) for all
# source/target cluster combinations (ECs, FAPs, MuSCs, Resolving-Mac), adding the
# previously-missing "Resolving-Mac" source-cluster rows and refreshing existing values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema4d"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.90715"
$ws.Range("H2").Value = [double]"2.72145"
$ws.Range("I2").Value = [double]"0.01717809939998381"
$ws.Range("J2").Value = [double]"0.01717809939998381"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"2.885873333333334"
$ws.Range("N2").Value = [double]"8.657620000000001"
$ws.Range("O2").Value = [double]"0.3070415651026022"
$ws.Range("P2").Value = [double]"0.3070415651026022"
$ws.Range("Q2").Value = [double]"2.617919994333334"
$ws.Range("R2").Value = [double]"23.561279949"
$ws.Range("S2").Value = [double]"0.005274390525259102"
$ws.Range("T2").Value = [double]"0.0052743905252591"

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema4d"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"0.90715"
$ws.Range("H3").Value = [double]"2.72145"
$ws.Range("I3").Value = [double]"0.01717809939998381"
$ws.Range("J3").Value = [double]"0.01717809939998381"
$ws.Range("K3").Value = [double]"3"
$ws.Range("L3").Value = [double]"1"
$ws.Range("M3").Value = [double]"3.165953666666667"
$ws.Range("N3").Value = [double]"9.497861"
$ws.Range("O3").Value = [double]"0.3368406220840099"
$ws.Range("P3").Value = [double]"0.3368406220840099"
$ws.Range("Q3").Value = [double]"2.871994868716667"
$ws.Range("R3").Value = [double]"25.84795381845"
$ws.Range("S3").Value = [double]"0.005786281688111505"
$ws.Range("T3").Value = [double]"0.005786281688111504"

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema4d"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"0.90715"
$ws.Range("H4").Value = [double]"2.72145"
$ws.Range("I4").Value = [double]"0.01717809939998381"
$ws.Range("J4").Value = [double]"0.01717809939998381"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"3.327024333333334"
$ws.Range("N4").Value = [double]"9.981073"
$ws.Range("O4").Value = [double]"0.3539776838580724"
$ws.Range("P4").Value = [double]"0.3539776838580724"
$ws.Range("Q4").Value = [double]"3.018110123983333"
$ws.Range("R4").Value = [double]"27.16299111585"
$ws.Range("S4").Value = [double]"0.006080663838690013"
$ws.Range("T4").Value = [double]"0.006080663838690012"

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sema4d"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"0.90715"
$ws.Range("H5").Value = [double]"2.72145"
$ws.Range("I5").Value = [double]"0.01717809939998381"
$ws.Range("J5").Value = [double]"0.01717809939998381"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.020115"
$ws.Range("N5").Value = [double]"0.060345"
$ws.Range("O5").Value = [double]"0.002140128955315263"
$ws.Range("P5").Value = [double]"0.002140128955315263"
$ws.Range("Q5").Value = [double]"0.01824732225"
$ws.Range("R5").Value = [double]"0.16422590025"
$ws.Range("S5").Value = [double]"3.67633479231891E-05"
$ws.Range("T5").Value = [double]"3.67633479231891E-05"

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema4d"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"1.252512666666667"
$ws.Range("H6").Value = [double]"3.757538"
$ws.Range("I6").Value = [double]"0.0237180037344858"
$ws.Range("J6").Value = [double]"0.0237180037344858"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"2.885873333333334"
$ws.Range("N6").Value = [double]"8.657620000000001"
$ws.Range("O6").Value = [double]"0.3070415651026022"
$ws.Range("P6").Value = [double]"0.3070415651026022"
$ws.Range("Q6").Value = [double]"3.614592904395557"
$ws.Range("R6").Value = [double]"32.53133613956001"
$ws.Range("S6").Value = [double]"0.007282412987745884"
$ws.Range("T6").Value = [double]"0.007282412987745884"

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema4d"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"1.252512666666667"
$ws.Range("H7").Value = [double]"3.757538"
$ws.Range("I7").Value = [double]"0.0237180037344858"
$ws.Range("J7").Value = [double]"0.0237180037344858"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"3.165953666666667"
$ws.Range("N7").Value = [double]"9.497861"
$ws.Range("O7").Value = [double]"0.3368406220840099"
$ws.Range("P7").Value = [double]"0.3368406220840099"
$ws.Range("Q7").Value = [double]"3.965397069579778"
$ws.Range("R7").Value = [double]"35.688573626218"
$ws.Range("S7").Value = [double]"0.007989187132515068"
$ws.Range("T7").Value = [double]"0.007989187132515067"

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sema4d"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"1.252512666666667"
$ws.Range("H8").Value = [double]"3.757538"
$ws.Range("I8").Value = [double]"0.0237180037344858"
$ws.Range("J8").Value = [double]"0.0237180037344858"
$ws.Range("K8").Value = [double]"3"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"3.327024333333334"
$ws.Range("N8").Value = [double]"9.981073"
$ws.Range("O8").Value = [double]"0.3539776838580724"
$ws.Range("P8").Value = [double]"0.3539776838580724"
$ws.Range("Q8").Value = [double]"4.167140119808223"
$ws.Range("R8").Value = [double]"37.50426107827401"
$ws.Range("S8").Value = [double]"0.008395644027670396"
$ws.Range("T8").Value = [double]"0.008395644027670394"

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sema4d"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"1.252512666666667"
$ws.Range("H9").Value = [double]"3.757538"
$ws.Range("I9").Value = [double]"0.0237180037344858"
$ws.Range("J9").Value = [double]"0.0237180037344858"
$ws.Range("K9").Value = [double]"1"
$ws.Range("L9").Value = [double]"0.3333333333333333"
$ws.Range("M9").Value = [double]"0.020115"
$ws.Range("N9").Value = [double]"0.060345"
$ws.Range("O9").Value = [double]"0.002140128955315263"
$ws.Range("P9").Value = [double]"0.002140128955315263"
$ws.Range("Q9").Value = [double]"0.02519429229"
$ws.Range("R9").Value = [double]"0.22674863061"
$ws.Range("S9").Value = [double]"5.07595865544486E-05"
$ws.Range("T9").Value = [double]"5.07595865544486E-05"

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Sema4d"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = [double]"3"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"2.247832333333333"
$ws.Range("H10").Value = [double]"6.743497"
$ws.Range("I10").Value = [double]"0.04256571378106988"
$ws.Range("J10").Value = [double]"0.04256571378106987"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"2.885873333333334"
$ws.Range("N10").Value = [double]"8.657620000000001"
$ws.Range("O10").Value = [double]"0.3070415651026022"
$ws.Range("P10").Value = [double]"0.3070415651026022"
$ws.Range("Q10").Value = [double]"6.486959388571112"
$ws.Range("R10").Value = [double]"58.38263449714"
$ws.Range("S10").Value = [double]"0.0130694433790491"
$ws.Range("T10").Value = [double]"0.0130694433790491"

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Sema4d"
$ws.Range("C11").Value = "Erbb2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = [double]"3"
$ws.Range("F11").Value = [double]"1"
$ws.Range("G11").Value = [double]"2.247832333333333"
$ws.Range("H11").Value = [double]"6.743497"
$ws.Range("I11").Value = [double]"0.04256571378106988"
$ws.Range("J11").Value = [double]"0.04256571378106987"
$ws.Range("K11").Value = [double]"3"
$ws.Range("L11").Value = [double]"1"
$ws.Range("M11").Value = [double]"3.165953666666667"
$ws.Range("N11").Value = [double]"9.497861"
$ws.Range("O11").Value = [double]"0.3368406220840099"
$ws.Range("P11").Value = [double]"0.3368406220840099"
$ws.Range("Q11").Value = [double]"7.116533017768555"
$ws.Range("R11").Value = [double]"64.048797159917"
$ws.Range("S11").Value = [double]"0.01433786150946549"
$ws.Range("T11").Value = [double]"0.01433786150946549"

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Sema4d"
$ws.Range("C12").Value = "Erbb2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = [double]"3"
$ws.Range("F12").Value = [double]"1"
$ws.Range("G12").Value = [double]"2.247832333333333"
$ws.Range("H12").Value = [double]"6.743497"
$ws.Range("I12").Value = [double]"0.04256571378106988"
$ws.Range("J12").Value = [double]"0.04256571378106987"
$ws.Range("K12").Value = [double]"3"
$ws.Range("L12").Value = [double]"1"
$ws.Range("M12").Value = [double]"3.327024333333334"
$ws.Range("N12").Value = [double]"9.981073"
$ws.Range("O12").Value = [double]"0.3539776838580724"
$ws.Range("P12").Value = [double]"0.3539776838580724"
$ws.Range("Q12").Value = [double]"7.478592870253444"
$ws.Range("R12").Value = [double]"67.307335832281"
$ws.Range("S12").Value = [double]"0.01506731277598875"
$ws.Range("T12").Value = [double]"0.01506731277598875"

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Sema4d"
$ws.Range("C13").Value = "Erbb2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = [double]"3"
$ws.Range("F13").Value = [double]"1"
$ws.Range("G13").Value = [double]"2.247832333333333"
$ws.Range("H13").Value = [double]"6.743497"
$ws.Range("I13").Value = [double]"0.04256571378106988"
$ws.Range("J13").Value = [double]"0.04256571378106987"
$ws.Range("K13").Value = [double]"1"
$ws.Range("L13").Value = [double]"0.3333333333333333"
$ws.Range("M13").Value = [double]"0.020115"
$ws.Range("N13").Value = [double]"0.060345"
$ws.Range("O13").Value = [double]"0.002140128955315263"
$ws.Range("P13").Value = [double]"0.002140128955315263"
$ws.Range("Q13").Value = [double]"0.045215147385"
$ws.Range("R13").Value = [double]"0.406936326465"
$ws.Range("S13").Value = [double]"9.109611656652958E-05"
$ws.Range("T13").Value = [double]"9.109611656652958E-05"

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Sema4d"
$ws.Range("C14").Value = "Erbb2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = [double]"3"
$ws.Range("F14").Value = [double]"1"
$ws.Range("G14").Value = [double]"48.40102466666667"
$ws.Range("H14").Value = [double]"145.203074"
$ws.Range("I14").Value = [double]"0.9165381830844606"
$ws.Range("J14").Value = [double]"0.9165381830844604"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"2.885873333333334"
$ws.Range("N14").Value = [double]"8.657620000000001"
$ws.Range("O14").Value = [double]"0.3070415651026022"
$ws.Range("P14").Value = [double]"0.3070415651026022"
$ws.Range("Q14").Value = [double]"139.6792263915423"
$ws.Range("R14").Value = [double]"1257.11303752388"
$ws.Range("S14").Value = [double]"0.2814153182105482"
$ws.Range("T14").Value = [double]"0.2814153182105481"

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Sema4d"
$ws.Range("C15").Value = "Erbb2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = [double]"3"
$ws.Range("F15").Value = [double]"1"
$ws.Range("G15").Value = [double]"48.40102466666667"
$ws.Range("H15").Value = [double]"145.203074"
$ws.Range("I15").Value = [double]"0.9165381830844606"
$ws.Range("J15").Value = [double]"0.9165381830844604"
$ws.Range("K15").Value = [double]"3"
$ws.Range("L15").Value = [double]"1"
$ws.Range("M15").Value = [double]"3.165953666666667"
$ws.Range("N15").Value = [double]"9.497861"
$ws.Range("O15").Value = [double]"0.3368406220840099"
$ws.Range("P15").Value = [double]"0.3368406220840099"
$ws.Range("Q15").Value = [double]"153.2354015138571"
$ws.Range("R15").Value = [double]"1379.118613624714"
$ws.Range("S15").Value = [double]"0.3087272917539179"
$ws.Range("T15").Value = [double]"0.3087272917539178"

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Sema4d"
$ws.Range("C16").Value = "Erbb2"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = [double]"3"
$ws.Range("F16").Value = [double]"1"
$ws.Range("G16").Value = [double]"48.40102466666667"
$ws.Range("H16").Value = [double]"145.203074"
$ws.Range("I16").Value = [double]"0.9165381830844606"
$ws.Range("J16").Value = [double]"0.9165381830844604"
$ws.Range("K16").Value = [double]"3"
$ws.Range("L16").Value = [double]"1"
$ws.Range("M16").Value = [double]"3.327024333333334"
$ws.Range("N16").Value = [double]"9.981073"
$ws.Range("O16").Value = [double]"0.3539776838580724"
$ws.Range("P16").Value = [double]"0.3539776838580724"
$ws.Range("Q16").Value = [double]"161.0313868242669"
$ws.Range("R16").Value = [double]"1449.282481418402"
$ws.Range("S16").Value = [double]"0.3244340632157233"
$ws.Range("T16").Value = [double]"0.3244340632157233"

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Sema4d"
$ws.Range("C17").Value = "Erbb2"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = [double]"3"
$ws.Range("F17").Value = [double]"1"
$ws.Range("G17").Value = [double]"48.40102466666667"
$ws.Range("H17").Value = [double]"145.203074"
$ws.Range("I17").Value = [double]"0.9165381830844606"
$ws.Range("J17").Value = [double]"0.9165381830844604"
$ws.Range("K17").Value = [double]"1"
$ws.Range("L17").Value = [double]"0.3333333333333333"
$ws.Range("M17").Value = [double]"0.020115"
$ws.Range("N17").Value = [double]"0.060345"
$ws.Range("O17").Value = [double]"0.002140128955315263"
$ws.Range("P17").Value = [double]"0.002140128955315263"
$ws.Range("Q17").Value = [double]"0.9735866111700001"
$ws.Range("R17").Value = [double]"8.762279500530001"
$ws.Range("S17").Value = [double]"0.001961509904271096"
$ws.Range("T17").Value = [double]"0.001961509904271096"
